# Add 2022-Q4 data:
#  - New "2022-Q4" worksheet inserted right after "总计", before "2022-Q3",
#    holding the per-fund holdings table for the new quarter.
#  - "总计" (summary) sheet gets a new row 2 for 2022-Q4, with every
#    following row shifted down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a cell as literal TEXT even when it "looks like"
# a number (e.g. "99.60"), without leaving behind a NumberFormat-driven style
# change. We build the literal in a scratch cell via a text formula (="...")
# then Copy/PasteSpecial-Values it into the destination - PasteSpecial keeps
# the string type but drops the formula, and never touches cell styles.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($ws, $addr, [string]$val)
    $escaped = $val.Replace('"', '""')
    $scratch = $ws.Range("ZZ1000")
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
    $scratch.ClearContents()
}

# ---------------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q3" sheet (position 2) and place the copy
#    right before it, so it lands in position 2 too (pushing the original
#    "2022-Q3" - and everything after it - one slot to the right). This
#    gives the new sheet identical formatting (header style, index-column
#    style, page margins, etc.) to the sheets it sits next to.
# ---------------------------------------------------------------------------
$sheetQ3 = $wb.Worksheets.Item(2)
$sheetQ3.Copy($sheetQ3)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2) Overwrite the duplicated sheet's data with the real 2022-Q4 numbers.
#    Row 1 (headers) is already correct from the copy.
# ---------------------------------------------------------------------------
$q4Data = @(
    @("159745", "国泰中证全指建筑材料ETF",   "7.45", "99.60", "2.56", "0.1907", 10),
    @("004856", "广发中证全指建筑材料指数A", "7.50", "94.40", "2.37", "0.1778", 10),
    @("004857", "广发中证全指建筑材料指数C", "4.88", "94.40", "2.37", "0.1157", 10),
    @("516750", "富国中证全指建筑材料ETF",   "1.02", "99.35", "2.48", "0.0253", 10),
    @("004192", "招商中证500指数增强A",      "0.69", "93.06", "1.39", "0.0096", 4),
    @("004193", "招商中证500指数增强C",      "0.57", "93.06", "1.39", "0.0079", 4),
    @("159787", "易方达中证全指建筑材料ETF", "0.22", "96.55", "2.41", "0.0053", 10),
    @("003366", "浙商汇金中证转型成长指数",   "0.06", "93.82", "1.21", "0.0007", 10)
)

$lastExistingRow = $newSheet.UsedRange.Rows.Count   # rows already present after the sheet-copy

for ($i = 0; $i -lt $q4Data.Count; $i++) {
    $row = $i + 2
    $rec = $q4Data[$i]

    if ($row -gt $lastExistingRow) {
        # This row didn't exist on the source "2022-Q3" sheet (it only had 7
        # funds / 8 rows); clone formatting from the row above it so the
        # index column keeps its style.
        $newSheet.Range("A" + ($row - 1) + ":H" + ($row - 1)).Copy($newSheet.Range("A$row" + ":H$row"))
    }

    $newSheet.Range("A$row").Value = $i

    Set-TextValue $newSheet "B$row" $rec[0]
    Set-TextValue $newSheet "C$row" $rec[1]
    Set-TextValue $newSheet "D$row" $rec[2]
    Set-TextValue $newSheet "E$row" $rec[3]
    Set-TextValue $newSheet "F$row" $rec[4]
    Set-TextValue $newSheet "G$row" $rec[5]

    $newSheet.Range("H$row").Value = $rec[6]
}

# ---------------------------------------------------------------------------
# 3) Update the "总计" (summary) sheet: shift rows 2-8 down to 3-9 (copying
#    values+formatting together so the index-column style follows), then
#    write the brand-new 2022-Q4 summary row into row 2.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

for ($r = 8; $r -ge 2; $r--) {
    $src = $totalSheet.Range("A" + $r + ":D" + $r)
    $dst = $totalSheet.Range("A" + ($r + 1) + ":D" + ($r + 1))
    $src.Copy($dst)
}

$totalSheet.Range("A2").Value = 0
Set-TextValue $totalSheet "B2" "2022-Q4"
$totalSheet.Range("C2").Value = 8
$totalSheet.Range("D2").Value = 0.53

Write-Host "2022-Q4 sheet added and summary sheet updated"
